$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q3, pushing the
#    existing quarters (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3) down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.4
$summary.Range("A2").Value = 0

# The freshly-inserted row copied formatting from the row above (the
# header) for B2:D2 and left A2 unstyled - fix both so the row matches
# the look of the other data rows (only column A carries the bold style).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A is a plain 0-based row index (row number - 2); re-number it
# for every data row now that one extra row exists.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 2) Add a new "2022-Q3" sheet, right after "总计", by duplicating the
#    "2022-Q2" sheet (this carries over all sheet/column/cell formatting)
#    and then overwriting its values with the Q3 figures.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Columns B-G hold text-formatted figures (fund code / name / metrics);
# set them to Text format first so values like "000593" or "93.04" are
# not silently reinterpreted as numbers.
$textRange = $q3.Range("B2:G4")
$textRange.NumberFormat = "@"

$q3.Range("B2").Value = "118002"
$q3.Range("C2").Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$q3.Range("D2").Value = "1.85"
$q3.Range("E2").Value = "93.04"
$q3.Range("F2").Value = "7.17"
$q3.Range("G2").Value = "0.1326"
$q3.Range("H2").Value = 5

$q3.Range("B3").Value = "000593"
$q3.Range("C3").Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$q3.Range("D3").Value = "1.85"
$q3.Range("E3").Value = "93.04"
$q3.Range("F3").Value = "7.17"
$q3.Range("G3").Value = "0.1326"
$q3.Range("H3").Value = 5

$q3.Range("B4").Value = "005676"
$q3.Range("C4").Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$q3.Range("D4").Value = "1.85"
$q3.Range("E4").Value = "93.04"
$q3.Range("F4").Value = "7.17"
$q3.Range("G4").Value = "0.1326"
$q3.Range("H4").Value = 5

# Drop the Text-format override again now that the values are committed,
# so the cells end up with the same plain (unstyled) look as the rest of
# the sheet.
$textRange.Style = "Normal"
